$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4624.143
$ws.Range("J116").Value = 5295
$ws.Range("L116").Value = 5295
$ws.Range("N116").Value = -12179

$ws.Range("H135").Value = 2285.7144
$ws.Range("I135").Value = 1926.2
$ws.Range("J135").Value = 3184.5
$ws.Range("K135").Value = 17335.8
$ws.Range("L135").Value = 28660.5
$ws.Range("M135").Value = -14800.8
$ws.Range("N135").Value = -33730.5

$ws.Range("H137").Value = 2425.158
$ws.Range("I137").Value = 1600
$ws.Range("J137").Value = 2906.5
$ws.Range("K137").Value = 4800
$ws.Range("L137").Value = 8719.5
$ws.Range("M137").Value = -2250
$ws.Range("N137").Value = -13819.5

$ws.Range("H141").Value = 3560.7778
$ws.Range("I141").Value = 3530.875
$ws.Range("K141").Value = 10592.625
$ws.Range("M141").Value = -5412.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H37").Value = 59999.5
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("M37").Value = -19727

$ws.Range("H61").Value = 3570.2856
$ws.Range("I61").Value = 3570.2856
$ws.Range("K61").Value = 3570.2856
$ws.Range("M61").Value = -3358.2856

$ws.Range("H122").Value = 1220.7142
$ws.Range("I122").Value = 1220.7142
$ws.Range("K122").Value = 3662.1426
$ws.Range("M122").Value = -1212.1426

$ws.Range("H136").Value = 3570.2856
$ws.Range("I136").Value = 3570.2856
$ws.Range("K136").Value = 10710.8568
$ws.Range("M136").Value = -8160.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H36").Value = 933
$ws.Range("I36").Value = 933
$ws.Range("K36").Value = 933
$ws.Range("M36").Value = -399

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

$ws.Range("H107").Value = 4845.1113
$ws.Range("I107").Value = 1565.6364
$ws.Range("K107").Value = 1565.6364
$ws.Range("M107").Value = 354.3635999999999

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 187445
$ws.Range("J52").Value = 187445
$ws.Range("L52").Value = 187445
$ws.Range("N52").Value = -188033

$ws.Range("H58").Value = 2857.6
$ws.Range("J58").Value = 5747.5
$ws.Range("L58").Value = 5747.5
$ws.Range("N58").Value = -6153.5

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H136").Value = 2857.6
$ws.Range("J136").Value = 5747.5
$ws.Range("L136").Value = 17242.5
$ws.Range("N136").Value = -22342.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2031.75
$ws.Range("I5").Value = 1653.5555
$ws.Range("J5").Value = 3166.3333
$ws.Range("K5").Value = 4960.666499999999
$ws.Range("L5").Value = 9498.999899999999
$ws.Range("M5").Value = -4848.666499999999
$ws.Range("N5").Value = -9722.999899999999

$ws.Range("H7").Value = 89.5
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 97.40000000000001
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 292.2
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -516.2

$ws.Range("H15").Value = 174.83333
$ws.Range("J15").Value = 190
$ws.Range("L15").Value = 570
$ws.Range("N15").Value = -850

$ws.Range("H23").Value = 57.42857
$ws.Range("J23").Value = 90.666664
$ws.Range("L23").Value = 271.999992
$ws.Range("N23").Value = -741.999992

$ws.Range("H34").Value = 926.13336
$ws.Range("I34").Value = 169.5
$ws.Range("J34").Value = 1201.2727
$ws.Range("K34").Value = 508.5
$ws.Range("L34").Value = 3603.8181
$ws.Range("M34").Value = -424.5
$ws.Range("N34").Value = -3771.8181

$ws.Range("H42").Value = 4833.3335
$ws.Range("J42").Value = 5750
$ws.Range("L42").Value = 17250
$ws.Range("N42").Value = -18318

$ws.Range("H49").Value = 2997.5
$ws.Range("I49").Value = 2995
$ws.Range("J49").Value = 3000
$ws.Range("K49").Value = 8985
$ws.Range("L49").Value = 9000
$ws.Range("M49").Value = -8829
$ws.Range("N49").Value = -9312

$ws.Range("H107").Value = 478.1143
$ws.Range("I107").Value = 240.47058
$ws.Range("J107").Value = 702.55554
$ws.Range("K107").Value = 721.41174
$ws.Range("L107").Value = 2107.66662
$ws.Range("M107").Value = 1198.58826
$ws.Range("N107").Value = -5947.66662

$ws.Range("H132").Value = 5350.8
$ws.Range("J132").Value = 5501.25
$ws.Range("L132").Value = 49511.25
$ws.Range("N132").Value = -54571.25

$ws.Range("H134").Value = 3032.5
$ws.Range("I134").Value = 3032.5
$ws.Range("K134").Value = 9097.5
$ws.Range("M134").Value = -4027.5

$ws.Range("H135").Value = 2031.75
$ws.Range("I135").Value = 1653.5555
$ws.Range("J135").Value = 3166.3333
$ws.Range("K135").Value = 14881.9995
$ws.Range("L135").Value = 28496.9997
$ws.Range("M135").Value = -12346.9995
$ws.Range("N135").Value = -33566.9997

$ws.Range("H136").Value = 10990
$ws.Range("I136").Value = 8980
$ws.Range("K136").Value = 26940
$ws.Range("M136").Value = -21840

$ws.Range("H137").Value = 4581.75
$ws.Range("J137").Value = 5765.6665
$ws.Range("L137").Value = 17296.9995
$ws.Range("N137").Value = -27496.9995

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

$ws.Range("H113").Value = 7345.385
$ws.Range("I113").Value = 5772.75
$ws.Range("J113").Value = 8044.3335
$ws.Range("K113").Value = 5772.75
$ws.Range("L113").Value = 8044.3335
$ws.Range("M113").Value = -3602.75
$ws.Range("N113").Value = -12384.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1274.2
$ws.Range("I22").Value = 1299.75
$ws.Range("J22").Value = 1257.1666
$ws.Range("K22").Value = 1299.75
$ws.Range("L22").Value = 1257.1666
$ws.Range("M22").Value = -1004.75
$ws.Range("N22").Value = -1847.1666

$ws.Range("H27").Value = 1274.2
$ws.Range("I27").Value = 1299.75
$ws.Range("J27").Value = 1257.1666
$ws.Range("K27").Value = 1299.75
$ws.Range("L27").Value = 1257.1666
$ws.Range("M27").Value = -1192.75
$ws.Range("N27").Value = -1471.1666

$ws.Range("H61").Value = 6098.778
$ws.Range("I61").Value = 4677.8
$ws.Range("K61").Value = 4677.8
$ws.Range("M61").Value = -4475.8

$ws.Range("H74").Value = 47196.5
$ws.Range("I74").Value = 47196.5
$ws.Range("K74").Value = 47196.5
$ws.Range("M74").Value = -46198.5

$ws.Range("H77").Value = 47196.5
$ws.Range("I77").Value = 47196.5
$ws.Range("K77").Value = 141589.5
$ws.Range("M77").Value = -136597.5

$ws.Range("H80").Value = 24000
$ws.Range("J80").Value = 24000
$ws.Range("L80").Value = 24000
$ws.Range("N80").Value = -26246

$ws.Range("H83").Value = 24000
$ws.Range("J83").Value = 24000
$ws.Range("L83").Value = 72000
$ws.Range("N83").Value = -83232

$ws.Range("H113").Value = 6098.778
$ws.Range("I113").Value = 4677.8
$ws.Range("K113").Value = 4677.8
$ws.Range("M113").Value = -2507.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1791.8334
$ws.Range("I100").Value = 1791.8334
$ws.Range("K100").Value = 3583.6668
$ws.Range("M100").Value = -3042.6668

$ws.Range("H107").Value = 958
$ws.Range("I107").Value = 972
$ws.Range("J107").Value = 944
$ws.Range("K107").Value = 2916
$ws.Range("L107").Value = 2832
$ws.Range("M107").Value = -996
$ws.Range("N107").Value = -6672

Write-Output "done"